$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 363; this shifts the old rows 363-377 down to 365-379
$ws.Rows.Item(363).Resize(2).Insert()

# Populate new row 363 (Primera, Region del Maule)
$ws.Cells.Item(363,1).Value = 7
$ws.Cells.Item(363,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(363,3).Value = "Ñuble"
$ws.Cells.Item(363,4).Value = 44939
$ws.Cells.Item(363,5).Value = 16
$ws.Cells.Item(363,6).Value = 100112008
$ws.Cells.Item(363,7).Value = "Coliflor"
$ws.Cells.Item(363,8).Value = "Sin especificar"
$ws.Cells.Item(363,9).Value = "Primera"
$ws.Cells.Item(363,10).Value = 100
$ws.Cells.Item(363,11).Value = 800
$ws.Cells.Item(363,12).Value = 800
$ws.Cells.Item(363,13).Value = 800
$ws.Cells.Item(363,14).Value = "$/unidad"
$ws.Cells.Item(363,15).Value = "Región del Maule"
$ws.Cells.Item(363,16).Value = 800
$ws.Cells.Item(363,17).Value = 1
$ws.Cells.Item(363,18).Value = "Hortaliza"

# Populate new row 364 (Segunda, Region del Maule)
$ws.Cells.Item(364,1).Value = 7
$ws.Cells.Item(364,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(364,3).Value = "Ñuble"
$ws.Cells.Item(364,4).Value = 44939
$ws.Cells.Item(364,5).Value = 16
$ws.Cells.Item(364,6).Value = 100112008
$ws.Cells.Item(364,7).Value = "Coliflor"
$ws.Cells.Item(364,8).Value = "Sin especificar"
$ws.Cells.Item(364,9).Value = "Segunda"
$ws.Cells.Item(364,10).Value = 100
$ws.Cells.Item(364,11).Value = 700
$ws.Cells.Item(364,12).Value = 700
$ws.Cells.Item(364,13).Value = 700
$ws.Cells.Item(364,14).Value = "$/unidad"
$ws.Cells.Item(364,15).Value = "Región del Maule"
$ws.Cells.Item(364,16).Value = 700
$ws.Cells.Item(364,17).Value = 1
$ws.Cells.Item(364,18).Value = "Hortaliza"

# Ensure the date cells keep the same date style as the rest of column D (style index used elsewhere)
$ws.Cells.Item(364,4).NumberFormat = $ws.Cells.Item(365,4).NumberFormat
$ws.Cells.Item(363,4).NumberFormat = $ws.Cells.Item(365,4).NumberFormat
